# Adds a "2022-Q1" fund-holdings sheet (positioned right before the "总计"
# summary sheet) and updates the "总计" summary sheet with a new row for it.
#
# Strategy (chosen so the resulting sheetId / file layout matches how Excel
# itself would produce this edit):
#   1. The existing "总计" sheet is renamed to "2022-Q1" and its old
#      (4-column) content is cleared and replaced by the new 8-fund holdings
#      table (columns A-H).
#   2. A fresh copy of that sheet is made and placed right after it, renamed
#      back to "总计"; its content is cleared and rebuilt as the updated
#      summary table (old rows + the new 2022-Q1 row on top).
#
# xlPasteFormats is used (copy a cell that already carries the workbook's
# "index / header" style, then PasteSpecial just the formatting) so that the
# existing style (bold + thin border + centered) is reused instead of a new,
# slightly different style being created.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# A style-"2" cell we can borrow formatting from (bold+border+centered).
$styleDonor = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------------
# Step 1: repurpose the current "总计" sheet into the new "2022-Q1" sheet
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"
$q1.UsedRange.ClearContents()

# Header row (B1:H1) styling, copied from an existing header cell.
$styleDonor.Range("G1").Copy()
$q1.Range("B1:H1").PasteSpecial($xlPasteFormats)

# Index column (A2:A9) styling, copied from an existing index cell.
$styleDonor.Range("A2").Copy()
$q1.Range("A2:A9").PasteSpecial($xlPasteFormats)

# Header labels.
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

function Set-FundRow {
    param($sheet, $row, $code, $name, $scale, $stockPos, $posRatio, $marketValue, $rank)

    $sheet.Cells.Item($row, 1).Value = ($row - 2)

    # Columns B, D, E, F, G hold values that look numeric (fund codes with
    # leading zeros, or decimals whose trailing/leading zeros must survive)
    # so they are forced to Text before assignment, exactly as in the
    # source data.
    $cB = $sheet.Cells.Item($row, 2)
    $cB.NumberFormat = "@"
    $cB.Value = $code

    $sheet.Cells.Item($row, 3).Value = $name

    $cD = $sheet.Cells.Item($row, 4)
    $cD.NumberFormat = "@"
    $cD.Value = $scale

    $cE = $sheet.Cells.Item($row, 5)
    $cE.NumberFormat = "@"
    $cE.Value = $stockPos

    $cF = $sheet.Cells.Item($row, 6)
    $cF.NumberFormat = "@"
    $cF.Value = $posRatio

    $cG = $sheet.Cells.Item($row, 7)
    $cG.NumberFormat = "@"
    $cG.Value = $marketValue

    $sheet.Cells.Item($row, 8).Value = $rank
}

Set-FundRow $q1 2 "006567" "中泰星元价值优选灵活配置混合"           "44.13" "81.95" "4.42" "1.9505" 8
Set-FundRow $q1 3 "013776" "中泰兴为价值精选混合A"                 "20.31" "85.34" "4.18" "0.8490" 9
Set-FundRow $q1 4 "006624" "中泰玉衡价值优选混合"                  "17.75" "81.95" "4.36" "0.7739" 8
Set-FundRow $q1 5 "013777" "中泰兴为价值精选混合C"                 "8.71"  "85.34" "4.18" "0.3641" 9
Set-FundRow $q1 6 "007592" "华夏价值精选混合"                      "2.55"  "94.58" "4.14" "0.1056" 10
Set-FundRow $q1 7 "512190" "浙商汇金中证浙江凤凰行动50ETF"         "0.51"  "98.94" "6.54" "0.0334" 4
Set-FundRow $q1 8 "930602" "国信价值智选混合型集合资产管理计划"    "0.50"  "67.38" "4.24" "0.0212" 7
Set-FundRow $q1 9 "516570" "易方达中证石化产业交易型开放式指数证券投资基金" "0.36" "96.03" "2.94" "0.0106" 10

# ---------------------------------------------------------------------
# Step 2: create the new "总计" sheet as a copy of "2022-Q1" placed right
# after it, then rebuild it as the updated summary table.
# ---------------------------------------------------------------------
$q1.Copy($null, $q1)
$total = $wb.Worksheets.Item("2022-Q1 (2)")
$total.Name = "总计"
$total.UsedRange.ClearContents()

$styleDonor.Range("G1").Copy()
$total.Range("B1:D1").PasteSpecial($xlPasteFormats)

$styleDonor.Range("A2").Copy()
$total.Range("A2:A7").PasteSpecial($xlPasteFormats)

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

function Set-TotalRow {
    param($sheet, $row, $label, $count, $value)
    $sheet.Cells.Item($row, 1).Value = ($row - 2)
    $sheet.Cells.Item($row, 2).Value = $label
    $sheet.Cells.Item($row, 3).Value = $count
    $sheet.Cells.Item($row, 4).Value = $value
}

Set-TotalRow $total 2 "2022-Q1" 8  4.11
Set-TotalRow $total 3 "2021-Q4" 9  3.85
Set-TotalRow $total 4 "2021-Q3" 11 4.15
Set-TotalRow $total 5 "2021-Q2" 12 18.21
Set-TotalRow $total 6 "2021-Q1" 19 20.15
Set-TotalRow $total 7 "2020-Q4" 12 7.3

Write-Output "Sheets now: $($wb.Worksheets.Count)"
foreach ($s in $wb.Worksheets) {
    Write-Output (" - " + $s.Name)
}
